# Update report title
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C1").Value = "Informe Reclamacion Por Cliente a Corte: 10 - junio - 2022"

# Row 3: keep only "No." (1) and the order-type text; clear the rest
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = "77 - SALIDA ARECLAMACION OFERTA NR -005"
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()

# Row 4
$ws.Range("C4").Value = "SO5177"
$ws.Range("D4").Value = "Z7 - SALIDA RECLAMACION PROMOCIONAL/ OFERTAS PORTOS 008"
$ws.Range("E4").Value = "BODEGA-040-DIGITAL"
$ws.Range("H4").Value = "CAMILA DELGADO"
$ws.Range("I4").Value = "2022-03-29 13:59:23"
$ws.Range("J4").Value = "cancel"

# Row 5
$ws.Range("C5").Value = "SO4965"
$ws.Range("D5").Value = "Z7 - SALIDA RECLAMACION PROMOCIONAL/ OFERTAS PORTOS 008"
$ws.Range("E5").Value = "BODEGA-040-DIGITAL"
$ws.Range("H5").Value = "DUITAMA MARIA DE JESUS"
$ws.Range("I5").Value = "2022-03-25 17:16:29"
$ws.Range("J5").Value = "cancel"

# Row 6
$ws.Range("C6").Value = "SO5383"
$ws.Range("D6").Value = "Z7 - SALIDA RECLAMACION PROMOCIONAL/ OFERTAS PORTOS 008"
$ws.Range("E6").Value = "BODEGA-040-DIGITAL"
$ws.Range("H6").Value = "RICHARD ALEJANDRO ALVAREZ"
$ws.Range("I6").Value = "2022-04-08 16:07:40"
$ws.Range("J6").Value = "cancel"

# Row 7
$ws.Range("C7").Value = "SO5581"
$ws.Range("D7").Value = "Z7 - SALIDA RECLAMACION PROMOCIONAL/ OFERTAS PORTOS 008"
$ws.Range("E7").Value = "BODEGA-040-DIGITAL"
$ws.Range("H7").Value = "MACIAS TAMAYO HECTOR ARIEL"
$ws.Range("I7").Value = "2022-04-21 13:54:03"
$ws.Range("J7").Value = "cancel"

# Row 8 (new data row, was blank before)
$ws.Range("A8").Value = 6
$ws.Range("C8").Value = "SO5861"
$ws.Range("D8").Value = "Z7 - SALIDA RECLAMACION PROMOCIONAL/ OFERTAS PORTOS 008"
$ws.Range("E8").Value = "BODEGA-040-DIGITAL"
$ws.Range("H8").Value = "JAIR GARCIA"
$ws.Range("I8").Value = "2022-04-25 21:42:16"
$ws.Range("J8").Value = "cancel"

# Row 9 (new data row, was blank before)
$ws.Range("A9").Value = 7
$ws.Range("C9").Value = "SO5341"
$ws.Range("D9").Value = "Z7 - SALIDA RECLAMACION PROMOCIONAL/ OFERTAS PORTOS 008"
$ws.Range("E9").Value = "BODEGA-040-DIGITAL"
$ws.Range("H9").Value = "MACIAS TAMAYO HECTOR ARIEL"
$ws.Range("I9").Value = "2022-03-30 15:20:59"
$ws.Range("J9").Value = "cancel"

# Extend the sheet dimension down to row 22 by inserting rows before the
# trailing blank rows (10-15), which shifts them to 17-22 and grows the
# used range without introducing any new cell styles.
$ws.Range("A10:A16").EntireRow.Insert()
